# Hortaliza, Vega Monumental Concepción - Berenjena
# A new weekly price observation was recorded; insert it as a new row 23
# (pushing the existing row 23..81 data down to 24..82, which is exactly
# what the underlying CSV/export's "insert a new reading at the top of the
# series" edit produces).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23:81 down to 24:82, inheriting row 23's formatting
# (keeps the D-column date style `s=2` on the new row).
$ws.Rows(23).Insert()

# Populate the newly inserted row 23 with the new observation.
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44659
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112001
$ws.Range("G23").Value = "Berenjena"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 9000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 9500
$ws.Range("N23").Value = "$/caja 60 unidades"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 158
$ws.Range("Q23").Value = 60
$ws.Range("R23").Value = "Hortaliza"
